# Apply content updates for 2024-08-28 Wednesday worksheet
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-27 Tuesday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-08-28 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("44÷3=14, 2", $false, $false, $false, $false, $false, $true, 1, $false, "80÷8=10, 0", 2) | Out-Null
$d.Content.Find.Execute("10÷3=3, 1", $false, $false, $false, $false, $false, $true, 1, $false, "89÷4=22, 1", 2) | Out-Null
$d.Content.Find.Execute("22÷3=7, 1", $false, $false, $false, $false, $false, $true, 1, $false, "56÷9=6, 2", 2) | Out-Null
$d.Content.Find.Execute("28÷5=5, 3", $false, $false, $false, $false, $false, $true, 1, $false, "96÷7=13, 5", 2) | Out-Null
$d.Content.Find.Execute("55÷5=11, 0", $false, $false, $false, $false, $false, $true, 1, $false, "87÷8=10, 7", 2) | Out-Null
$d.Content.Find.Execute("66÷6=11, 0", $false, $false, $false, $false, $false, $true, 1, $false, "39÷3=13, 0", 2) | Out-Null
$d.Content.Find.Execute("45÷5=9, 0", $false, $false, $false, $false, $false, $true, 1, $false, "66÷6=11, 0", 2) | Out-Null
$d.Content.Find.Execute("72÷7=10, 2", $false, $false, $false, $false, $false, $true, 1, $false, "44÷4=11, 0", 2) | Out-Null
$d.Content.Find.Execute("95÷7=13, 4", $false, $false, $false, $false, $false, $true, 1, $false, "55÷6=9, 1", 2) | Out-Null
$d.Content.Find.Execute("51÷4=12, 3", $false, $false, $false, $false, $false, $true, 1, $false, "96÷8=12, 0", 2) | Out-Null
$d.Content.Find.Execute("29÷4=7, 1", $false, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$d.Content.Find.Execute("83÷6=13, 5", $false, $false, $false, $false, $false, $true, 1, $false, "52÷2=26, 0", 2) | Out-Null
$d.Content.Find.Execute("70÷7=10, 0", $false, $false, $false, $false, $false, $true, 1, $false, "68÷8=8, 4", 2) | Out-Null
$d.Content.Find.Execute("81÷5=16, 1", $false, $false, $false, $false, $false, $true, 1, $false, "15÷8=1, 7", 2) | Out-Null
$d.Content.Find.Execute("11÷2=5, 1", $false, $false, $false, $false, $false, $true, 1, $false, "29÷9=3, 2", 2) | Out-Null
$d.Content.Find.Execute("90÷4=22, 2", $false, $false, $false, $false, $false, $true, 1, $false, "27÷3=9, 0", 2) | Out-Null
$d.Content.Find.Execute("29÷6=4, 5", $false, $false, $false, $false, $false, $true, 1, $false, "24÷5=4, 4", 2) | Out-Null
$d.Content.Find.Execute("16÷6=2, 4", $false, $false, $false, $false, $false, $true, 1, $false, "42÷5=8, 2", 2) | Out-Null
$d.Content.Find.Execute("46÷8=5, 6", $false, $false, $false, $false, $false, $true, 1, $false, "40÷5=8, 0", 2) | Out-Null
$d.Content.Find.Execute("20÷2=10, 0", $false, $false, $false, $false, $false, $true, 1, $false, "25÷3=8, 1", 2) | Out-Null
$d.Content.Find.Execute("90÷6=15, 0", $false, $false, $false, $false, $false, $true, 1, $false, "21÷7=3, 0", 2) | Out-Null
$d.Content.Find.Execute("28÷8=3, 4", $false, $false, $false, $false, $false, $true, 1, $false, "30÷5=6, 0", 2) | Out-Null
$d.Content.Find.Execute("43÷4=10, 3", $false, $false, $false, $false, $false, $true, 1, $false, "27÷2=13, 1", 2) | Out-Null
$d.Content.Find.Execute("94÷2=47, 0", $false, $false, $false, $false, $false, $true, 1, $false, "37÷3=12, 1", 2) | Out-Null
$d.Content.Find.Execute("89÷8=11, 1", $false, $false, $false, $false, $false, $true, 1, $false, "14÷9=1, 5", 2) | Out-Null
